$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:A31").Value = 1
$ws.Range("A3").Select() | Out-Null
